$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 98

$kRange = $ws.Range("K2:K$lastRow")
$kRange.Formula = "=E2/D2"

$lRange = $ws.Range("L2:L$lastRow")
$lRange.Formula = "=H2/F2"

$ws.Calculate()

$ws.Range("N5").Select() | Out-Null
